$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "Polychronidou P"
$ws.Range("B10").Value = "Kydros D"
$ws.Range("C10").Value = "tralala"
$ws.Range("D10").Value = 332123
$ws.Range("E10").Value = 2022
$ws.Range("F10").Value = "statistics, econometrics"

# Row 11
$ws.Range("A11").Value = "Magulios G"
$ws.Range("B11").Value = "Polychronidou P"
$ws.Range("C11").Value = "kikiki"
$ws.Range("D11").Value = 11111
$ws.Range("E11").Value = 2023
$ws.Range("F11").Value = "economics"

# Row 12
$ws.Range("A12").Value = "Polychronidou P"
$ws.Range("B12").Value = "Vrana V"
$ws.Range("C12").Value = "lalala"
$ws.Range("D12").Value = 22222
$ws.Range("E12").Value = 2022
$ws.Range("F12").Value = "operations research"

# Update the active selection to B10 (matches the edited sheetView selection)
$ws.Range("B10").Select()
